# RUNMANAGER sheet: the "Module" and "Product" columns (F and G) were
# transposed - Product (with its header + values) now comes before Module.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RUNMANAGER")

$lastRow = 6
for ($r = 1; $r -le $lastRow; $r++) {
    $colF = $ws.Cells.Item($r, 6).Value2
    $colG = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 6).Value = $colG
    $ws.Cells.Item($r, 7).Value = $colF
}

# Selection moved from I9 to H9.
$null = $ws.Range("H9").Select()
